$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 9 (week 8) assignments
$ws.Range("B9").Value = "Harun"
$ws.Range("C9").Value = "Mohamed"

# Update selected cell to D9
$ws.Range("D9").Select()
